# natmiOut/YoungD7/LR-pairs_lrc2p/Serpinf1-Plxdc1.xlsx
# "Natmi following Dr Hou advice"
#
# The ligand/receptor interaction table is recomputed: every Sending-cluster x
# Target-cluster combination across the four clusters (ECs, FAPs, M2, sCs) is
# now included (16 rows instead of 12), and every numeric metric is refreshed
# with the newly recalculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-17 (columns A-T), in final/target form.
$data = @(
    @("ECs","Serpinf1","Plxdc1","ECs",1,0.3333333333333333,1.026875,3.080625,0.001791230048025246,0.001791230048025246,3,1,2.181962,6.545886,0.1718849556529879,0.1718849556529879,2.24060222875,20.16542005875,0.0003078854973691188,0.0003078854973691189),
    @("ECs","Serpinf1","Plxdc1","FAPs",1,0.3333333333333333,1.026875,3.080625,0.001791230048025246,0.001791230048025246,3,1,1.926782,5.780346,0.1517830460030813,0.1517830460030813,1.97856426625,17.80707839625,0.0002718783527815175,0.0002718783527815175),
    @("ECs","Serpinf1","Plxdc1","M2",1,0.3333333333333333,1.026875,3.080625,0.001791230048025246,0.001791230048025246,3,1,2.851342333333334,8.554027000000001,0.2246156672373245,0.2246156672373245,2.927972158541667,26.351749426875,0.0004023383324127355,0.0004023383324127355),
    @("ECs","Serpinf1","Plxdc1","sCs",1,0.3333333333333333,1.026875,3.080625,0.001791230048025246,0.001791230048025246,3,1,5.73423,17.20269,0.4517163311066063,0.4517163311066063,5.88833743125,52.99503688125,0.0008091278654618743,0.0008091278654618743),
    @("FAPs","Serpinf1","Plxdc1","ECs",3,1,543.3469646666667,1630.040894,0.9477876173967086,0.9477876173967087,3,1,2.181962,6.545886,0.1718849556529879,0.1718849556529879,1185.56242971801,10670.06186746209,0.1629104325846843,0.1629104325846844),
    @("FAPs","Serpinf1","Plxdc1","FAPs",3,1,543.3469646666667,1630.040894,0.9477876173967086,0.9477876173967087,3,1,1.926782,5.780346,0.1517830460030813,0.1517830460030813,1046.91115127437,9422.200361469324,0.1438580915324754,0.1438580915324755),
    @("FAPs","Serpinf1","Plxdc1","M2",3,1,543.3469646666667,1630.040894,0.9477876173967086,0.9477876173967087,3,1,2.851342333333334,8.554027000000001,0.2246156672373245,0.2246156672373245,1549.268202042238,13943.41381838014,0.2128879480808357,0.2128879480808358),
    @("FAPs","Serpinf1","Plxdc1","sCs",3,1,543.3469646666667,1630.040894,0.9477876173967086,0.9477876173967087,3,1,5.73423,17.20269,0.4517163311066063,0.4517163311066063,3115.67646520054,28041.08818680486,0.4281311451987131,0.4281311451987131),
    @("M2","Serpinf1","Plxdc1","ECs",3,1,0.9942943333333334,2.982883,0.001734397941763016,0.001734397941763016,3,1,2.181962,6.545886,0.1718849556529879,0.1718849556529879,2.169512452148667,19.525612069338,0.0002981169133045695,0.0002981169133045695),
    @("M2","Serpinf1","Plxdc1","FAPs",3,1,0.9942943333333334,2.982883,0.001734397941763016,0.001734397941763016,3,1,1.926782,5.780346,0.1517830460030813,0.1517830460030813,1.915788424168667,17.242095817518,0.0002632522025822654,0.0002632522025822654),
    @("M2","Serpinf1","Plxdc1","M2",3,1,0.9942943333333334,2.982883,0.001734397941763016,0.001734397941763016,3,1,2.851342333333334,8.554027000000001,0.2246156672373245,0.2246156672373245,2.835073524426778,25.51566171984101,0.0003895729509441421,0.0003895729509441421),
    @("M2","Serpinf1","Plxdc1","sCs",3,1,0.9942943333333334,2.982883,0.001734397941763016,0.001734397941763016,3,1,5.73423,17.20269,0.4517163311066063,0.4517163311066063,5.701512395030001,51.31361155527001,0.0007834558749320389,0.0007834558749320389),
    @("sCs","Serpinf1","Plxdc1","ECs",3,1,27.91110566666667,83.733317,0.04868675461350316,0.04868675461350316,3,1,2.181962,6.545886,0.1718849556529879,0.1718849556529879,60.90097194265133,548.108747483862,0.008368520657629894,0.008368520657629896),
    @("sCs","Serpinf1","Plxdc1","FAPs",3,1,27.91110566666667,83.733317,0.04868675461350316,0.04868675461350316,3,1,1.926782,5.780346,0.1517830460030813,0.1517830460030813,53.77861599863134,484.007543987682,0.007389823915242081,0.007389823915242081),
    @("sCs","Serpinf1","Plxdc1","M2",3,1,27.91110566666667,83.733317,0.04868675461350316,0.04868675461350316,3,1,2.851342333333334,8.554027000000001,0.2246156672373245,0.2246156672373245,79.58411715750657,716.2570544175591,0.0109358078731319,0.0109358078731319),
    @("sCs","Serpinf1","Plxdc1","sCs",3,1,27.91110566666667,83.733317,0.04868675461350316,0.04868675461350316,3,1,5.73423,17.20269,0.4517163311066063,0.4517163311066063,160.04869944697,1440.43829502273,0.02199260216749929,0.02199260216749929)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $row = $startRow + $i
    for ($col = 1; $col -le $rowValues.Count; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowValues[$col - 1]
    }
}
